# Apply the edit described by the diff:
# - Rename sheet "Employee" -> "Employees"
# - Replace placeholder name values "Hola"/"Hola2"/"Hola3" -> "Uno"/"Dos"/"Tres"
# - Refresh the "Date Of Birth" timestamps in column C (same dates, new time-of-day)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "Employees"

# Update the Name values (column A, rows 2-4)
$ws.Range("A2").Value = "Uno"
$ws.Range("A3").Value = "Dos"
$ws.Range("A4").Value = "Tres"

# Update the Date Of Birth values (column C, rows 2-4); same calendar day, new time-of-day
$ws.Range("C2").Value = 33837.90597678241
$ws.Range("C3").Value = 24061.905976782407
$ws.Range("C4").Value = 31915.905976782407
